# This script reproduces a weekly data refresh: a new observation row is
# inserted right above the former row 78, pushing every following row
# down by one (old row 78 becomes 79, ..., old row 144 becomes 145), and
# the newly opened row 78 is filled in with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 - this shifts rows 78:144 down to 79:145
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with the new weekly record.
# Columns A, B, C (Mercado ID, Mercado, Región) stay consistent with the
# rest of the dataset for this sheet.
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = 44790
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 100112032
$ws.Range("G78").Value = "Zapallo italiano"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 250
$ws.Range("K78").Value = 20000
$ws.Range("L78").Value = 22000
$ws.Range("M78").Value = 20800
$ws.Range("N78").Value = "$/caja 60 unidades"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 347
$ws.Range("Q78").Value = 60
$ws.Range("R78").Value = "Hortaliza"
